$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03612509674085
$ws.Range("D2").Value = 1.038138645255278
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.04360766331726
$ws.Range("I2").Value = 1.034158167621293
$ws.Range("J2").Value = 1.041235120978439
$ws.Range("K2").Value = 1.04092753164607
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.046381051344862
$ws.Range("N2").Value = 1.042713794281605

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037922818510341
$ws.Range("D3").Value = 1.039509490123185
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.045453735789065
$ws.Range("I3").Value = 1.034605084211917
$ws.Range("J3").Value = 1.042672406023564
$ws.Range("K3").Value = 1.04210684621371
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.048035490544426
$ws.Range("N3").Value = 1.044153120436315

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039081994242461
$ws.Range("D4").Value = 1.040392815400964
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.046644571299179
$ws.Range("I4").Value = 1.03489117027549
$ws.Range("J4").Value = 1.043598137212408
$ws.Range("K4").Value = 1.042865721643377
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.049101911696669
$ws.Range("N4").Value = 1.04508016626963

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039568356608781
$ws.Range("D5").Value = 1.040763292616308
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.047144332772767
$ws.Range("I5").Value = 1.035010704669345
$ws.Range("J5").Value = 1.043986303652549
$ws.Range("K5").Value = 1.04318375555812
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.049549268613698
$ws.Range("N5").Value = 1.045468883950638

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039649963540637
$ws.Range("D6").Value = 1.040825446590854
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.047228194709199
$ws.Range("I6").Value = 1.035030732013508
$ws.Range("J6").Value = 1.044051419675552
$ws.Range("K6").Value = 1.0432370967421
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.049624325643864
$ws.Range("N6").Value = 1.045534092445863

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039088496766832
$ws.Range("D7").Value = 1.040397769146096
$ws.Range("E7").Value = 0.9943035907978915
$ws.Range("F7").Value = 1.046651252512191
$ws.Range("I7").Value = 1.034892770384848
$ws.Range("J7").Value = 1.043603327862383
$ws.Range("K7").Value = 1.042869975129233
$ws.Range("L7").Value = 0.9968970624459041
$ws.Range("M7").Value = 1.049107893075316
$ws.Range("N7").Value = 1.045085364290922

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036733499764159
$ws.Range("D8").Value = 1.038602703771005
$ws.Range("E8").Value = 0.9929600610674294
$ws.Range("F8").Value = 1.044232329252598
$ws.Range("I8").Value = 1.034309850256337
$ws.Range("J8").Value = 1.041721755612599
$ws.Range("K8").Value = 1.04132696813396
$ws.Range("L8").Value = 0.9958175282591053
$ws.Range("M8").Value = 1.046941038931198
$ws.Range("N8").Value = 1.043201119992724

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.032551644289019
$ws.Range("D9").Value = 1.035410585398478
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.039940694406426
$ws.Range("I9").Value = 1.033258679585634
$ws.Range("J9").Value = 1.038372632144506
$ws.Range("K9").Value = 1.03857508712197
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.04309047049321
$ws.Range("N9").Value = 1.039847240385157

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029740960714777
$ws.Range("D10").Value = 1.033262141114904
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.037058759538192
$ws.Range("I10").Value = 1.032541405614812
$ws.Range("J10").Value = 1.036116325308596
$ws.Range("K10").Value = 1.036717550341339
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.040500579549706
$ws.Range("N10").Value = 1.037587729334738

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028518224331002
$ws.Range("D11").Value = 1.032326813190843
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.03580563119657
$ws.Range("I11").Value = 1.032226824145286
$ws.Range("J11").Value = 1.03513350757703
$ws.Range("K11").Value = 1.035907587489115
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.03937345545901
$ws.Range("N11").Value = 1.036603515889263

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028063166494795
$ws.Range("D12").Value = 1.031978616776366
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.035339353743857
$ws.Range("I12").Value = 1.032109366906693
$ws.Range("J12").Value = 1.03476755113127
$ws.Range("K12").Value = 1.035605868150702
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.038953915341526
$ws.Range("N12").Value = 1.036237039743374

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028160818114856
$ws.Range("D13").Value = 1.032053341391661
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.035439408870862
$ws.Range("I13").Value = 1.032134589497593
$ws.Range("J13").Value = 1.034846090810584
$ws.Range("K13").Value = 1.035670627290563
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.039043948121127
$ws.Range("N13").Value = 1.036315690958041

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028480627181543
$ws.Range("D14").Value = 1.032298047055737
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.03576710521532
$ws.Range("I14").Value = 1.032217127531437
$ws.Range("J14").Value = 1.035103275857838
$ws.Range("K14").Value = 1.03588266497515
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.039338794141777
$ws.Range("N14").Value = 1.036573241237562

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028677554890438
$ws.Range("D15").Value = 1.032448715214166
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.035968901725187
$ws.Range("I15").Value = 1.032267901199631
$ws.Range("J15").Value = 1.035261616978118
$ws.Range("K15").Value = 1.036013193538212
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.039520341756653
$ws.Range("N15").Value = 1.036731807220395

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029821987742952
$ws.Range("D16").Value = 1.033324108157251
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.037141813291842
$ws.Range("I16").Value = 1.032562198528228
$ws.Range("J16").Value = 1.036181427364584
$ws.Range("K16").Value = 1.036771184766394
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.040575261428021
$ws.Range("N16").Value = 1.037652923843112

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030538319552758
$ws.Range("D17").Value = 1.033871858331491
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.037876132393959
$ws.Range("I17").Value = 1.032745728310911
$ws.Range("J17").Value = 1.036756827734192
$ws.Range("K17").Value = 1.037245132079132
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.041235447366518
$ws.Range("N17").Value = 1.038229141347245

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030955596232621
$ws.Range("D18").Value = 1.034190866879998
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.038303945161164
$ws.Range("I18").Value = 1.032852392953953
$ws.Range("J18").Value = 1.037091888628592
$ws.Range("K18").Value = 1.037521034725908
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.041619975715222
$ws.Range("N18").Value = 1.038564678066548

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031097784644029
$ws.Range("D19").Value = 1.034299558826049
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.038449733633239
$ws.Range("I19").Value = 1.032888697723374
$ws.Range("J19").Value = 1.037206041265883
$ws.Range("K19").Value = 1.037615018814966
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.041750997837291
$ws.Range("N19").Value = 1.038678992813676

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030461520736702
$ws.Range("D20").Value = 1.033813140196804
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.037797399117667
$ws.Range("I20").Value = 1.032726077194048
$ws.Range("J20").Value = 1.036695150825946
$ws.Range("K20").Value = 1.037194338278803
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.041164672394694
$ws.Range("N20").Value = 1.038167376850715

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028386475844045
$ws.Range("D21").Value = 1.032226008841492
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.035670629356579
$ws.Range("I21").Value = 1.032192838970432
$ws.Range("J21").Value = 1.035027566113266
$ws.Range("K21").Value = 1.035820249117798
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.039251993699323
$ws.Range("N21").Value = 1.036497423976468

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027076714440032
$ws.Range("D22").Value = 1.031223629682311
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.034328747557811
$ws.Range("I22").Value = 1.031854051499572
$ws.Range("J22").Value = 1.033973906109032
$ws.Range("K22").Value = 1.034951303089356
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.038044337013461
$ws.Range("N22").Value = 1.035442267654172

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027771535092079
$ws.Range("D23").Value = 1.031755440922609
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.035040557769296
$ws.Range("I23").Value = 1.032033985165704
$ws.Range("J23").Value = 1.034532969111224
$ws.Range("K23").Value = 1.035412427370062
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.038685027610695
$ws.Range("N23").Value = 1.036002124589952

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030496224507179
$ws.Range("D24").Value = 1.033839673892209
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.037832976850275
$ws.Range("I24").Value = 1.032734957877899
$ws.Range("J24").Value = 1.036723021697081
$ws.Range("K24").Value = 1.037217291487636
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.041196654247449
$ws.Range("N24").Value = 1.038195287301683

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.033636678012544
$ws.Range("D25").Value = 1.036239348915663
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.041053768237506
$ws.Range("I25").Value = 1.033533312533602
$ws.Range("J25").Value = 1.039242539011682
$ws.Range("K25").Value = 1.039290498090362
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.044089877445503
$ws.Range("N25").Value = 1.040718382619863
